# Corrections for exercise 9 (group member "cvx3958", column I/N block rows 4-6,
# 8-10 and 12-14 of Tabelle1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Achieved points for exercise 9 (column I) were corrected from 0 to 3 for the
# three "Achieved Points" rows.
$ws.Range("I5").Value = 3
$ws.Range("I9").Value = 3
$ws.Range("I13").Value = 3

# The "Max Points" column for exercise 9 (column N) is no longer graded, so the
# previously entered values are removed (cell formatting/style stays as-is).
$ws.Range("N4").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("N14").ClearContents()

# Leave the view on the cell that was active when the corrections were saved.
$ws.Activate()
$ws.Range("P9").Select()
